$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add new column J (mirrors column I formatting), update merged header ---
$ws.Range("I1").Copy($ws.Range("J1"))
$ws.Range("I2").Copy($ws.Range("J2"))
$ws.Range("J2").Value = "Heat exchanger network error [%]"

# Expand the C1:I1 merge to C1:J1. Merging re-derives per-cell borders (splitting
# the outline across the merged edge cells), so clear formatting first and then
# restore the original uniform header style (style index used by B1) afterwards.
$ws.Range("C1:I1").UnMerge()
$ws.Range("C1:J1").ClearFormats()
$ws.Range("C1:J1").Merge()
$ws.Range("B1").Copy()
$ws.Range("C1:J1").PasteSpecial(-4122, $false, $false, $false)

# --- 2. Rename metric headers in row 2 ---
$ws.Range("D2").Value = "Biodiesel production [MMGal/yr]"
$ws.Range("E2").Value = "Ethanol production [MMGal/yr]"
$ws.Range("F2").Value = "Electricity production [MMWhr/yr]"
$ws.Range("G2").Value = "Natural gas consumption [MMcf/yr]"
$ws.Range("H2").Value = "Productivity [MMGGE/yr]"

# --- 3. Reorder the lipidcane parameter labels (B4:B7) ---
$ws.Range("B4").Value = "Lipid content [dry wt. %]"
$ws.Range("B5").Value = "Lipid retention [%]"
$ws.Range("B6").Value = "Additional lipid extraction efficiency [%]"
$ws.Range("B7").Value = "Capacity [ton/hr]"

# --- 4. Update data values for rows 4-13, columns C,D,E,F,H,I, and new column J ---
$ws.Range("C4").Value = 0.01243446846537874
$ws.Range("D4").Value = 0.9656979317159171
$ws.Range("E4").Value = -0.9710174906486995
$ws.Range("F4").Value = 0.8874372113054882
$ws.Range("H4").Value = 0.2822626935945077
$ws.Range("I4").Value = 0.7359404309416171
$ws.Range("J4").Value = 0.8128816605944312

$ws.Range("C5").Value = 0.02509148922765957
$ws.Range("D5").Value = 0.05113234697329388
$ws.Range("E5").Value = 0.002162529782501191
$ws.Range("F5").Value = 0.02141216898448675
$ws.Range("H5").Value = 0.05909610255584409
$ws.Range("I5").Value = 0.0280392776335711
$ws.Range("J5").Value = 0.007775770778389963

$ws.Range("C6").Value = 0.04700840933633637
$ws.Range("D6").Value = 0.07846696883467874
$ws.Range("E6").Value = 0.006183976087359043
$ws.Range("F6").Value = -0.04582844180113766
$ws.Range("H6").Value = 0.0131443425097737
$ws.Range("I6").Value = -0.03494218690168747
$ws.Range("J6").Value = -0.00186048525922314

$ws.Range("C7").Value = 0.08224410194576406
$ws.Range("D7").Value = 0.1696790897951636
$ws.Range("E7").Value = 0.1878381166335246
$ws.Range("F7").Value = 0.3346338050333522
$ws.Range("H7").Value = 0.7480470345298813
$ws.Range("I7").Value = 0.651344777029791
$ws.Range("J7").Value = -0.05503422352114372

$ws.Range("C8").Value = 0.6458100750804029
$ws.Range("D8").Value = -0.00009880109195204365
$ws.Range("E8").Value = -0.01255370450214818
$ws.Range("F8").Value = -0.006136763861470554
$ws.Range("H8").Value = -0.02582272029690881
$ws.Range("I8").Value = -0.01007252075490083
$ws.Range("J8").Value = 0.005903371006685095

$ws.Range("C9").Value = 0.4012926813317072
$ws.Range("D9").Value = -0.008648459673938385
$ws.Range("E9").Value = 0.01109807113192284
$ws.Range("F9").Value = -0.004602154168086166
$ws.Range("H9").Value = 0.007134161373366454
$ws.Range("I9").Value = -0.00381772325670893
$ws.Range("J9").Value = 0.0004124588326765355

$ws.Range("C10").Value = 0.02957110303884412
$ws.Range("D10").Value = 0.004337107085484283
$ws.Range("E10").Value = -0.01196715311868612
$ws.Range("F10").Value = 0.005288213395528535
$ws.Range("H10").Value = -0.005754100742164029
$ws.Range("I10").Value = -0.008027295681091826
$ws.Range("J10").Value = -0.005842889612239713

$ws.Range("C11").Value = 0.1943246190689847
$ws.Range("D11").Value = 0.0001026954281078171
$ws.Range("E11").Value = -0.001599571359982854
$ws.Range("F11").Value = 0.0001663083906523356
$ws.Range("H11").Value = -0.004358397582335903
$ws.Range("I11").Value = 0.002187953655518146
$ws.Range("J11").Value = 0.001440911530258935

$ws.Range("C12").Value = 0.2058898845875954
$ws.Range("D12").Value = 0.1628570201142808
$ws.Range("E12").Value = 0.1141017462280698
$ws.Range("F12").Value = 0.2876288941291557
$ws.Range("H12").Value = 0.575809326232373
$ws.Range("I12").Value = 0.01952330132493205
$ws.Range("J12").Value = 0.0188931194838866

$ws.Range("C13").Value = -0.46037575111903
$ws.Range("D13").Value = 0.01137455632698225
$ws.Range("E13").Value = -0.01509476546779062
$ws.Range("F13").Value = 0.005795905863836234
$ws.Range("H13").Value = -0.009326392117055684
$ws.Range("I13").Value = 0.009080876331235051
$ws.Range("J13").Value = 0.005744321992254426

# --- 5. Remove the now-obsolete "Fermentation" row (row 14) ---
$ws.Range("A14:I14").EntireRow.Delete()
